# Update the "Training Dashboard" sheet with the new progress date (04-Nov-2025).
# For every data row (3 through 16):
#   - column H ("PERIOD TO EXPIRE") is decremented by 1 day
#   - column I ("LAST UPDATE") is updated from 03-Nov-2025 to 04-Nov-2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$newLastUpdate = "04-Nov-2025"

for ($row = 3; $row -le 16; $row++) {
    $periodCell = $ws.Range("H$row")
    $currentPeriod = $periodCell.Value()
    $periodCell.Value = $currentPeriod - 1

    # Prefix with an apostrophe so Excel stores the date-looking string as
    # literal text instead of auto-converting it into a date serial value.
    $ws.Range("I$row").Value = "'" + $newLastUpdate
}
